# Generate Report for Handback
#
# The source report re-uses a handful of shared strings across several
# cells (rows for 955bae9d-e21a-44d1-8b1f-05727881324f.md and
# cce7b96c-78f1-4b67-a336-2f5a6d986273.md on the "Overview" row, and the
# corresponding rows 4/5 on the "zh-cn"/"de-de" sheets, plus cross-sheet
# re-use of the very same timestamp/priority strings). A fresh report run
# bumps those shared values, so every cell that showed the old value now
# shows the new one.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# "2016-09-04 22:18:05" -> "2016-09-04 22:18:58"
# (shared by Overview!G4, Overview!G5, de-de!H4, de-de!H5)
$overview.Range("G4").Value = "2016-09-04 22:18:58"
$overview.Range("G5").Value = "2016-09-04 22:18:58"
$dede.Range("H4").Value = "2016-09-04 22:18:58"
$dede.Range("H5").Value = "2016-09-04 22:18:58"

# "ht" -> "mt"
# (shared by zh-cn!E4, zh-cn!E5, de-de!E4, de-de!E5)
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("E5").Value = "mt"
$dede.Range("E4").Value = "mt"
$dede.Range("E5").Value = "mt"

# "2016-09-04 22:17:54" -> "2016-09-04 22:18:53"
# (shared by zh-cn!H4, zh-cn!H5)
$zhcn.Range("H4").Value = "2016-09-04 22:18:53"
$zhcn.Range("H5").Value = "2016-09-04 22:18:53"

# "2016-09-04 22:18:27" -> "2016-09-04 22:19:16"
# (shared by zh-cn!K4, zh-cn!K5)
$zhcn.Range("K4").Value = "2016-09-04 22:19:16"
$zhcn.Range("K5").Value = "2016-09-04 22:19:16"

# "2016-09-04 22:18:34" -> "2016-09-04 22:19:23"
# (shared by de-de!K4, de-de!K5)
$dede.Range("K4").Value = "2016-09-04 22:19:23"
$dede.Range("K5").Value = "2016-09-04 22:19:23"
